$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 1.328826
    "H2" = 3.986478
    "I2" = 0.03246115949735631
    "J2" = 0.03246115949735631
    "M2" = 4.019613333333333
    "N2" = 12.05884
    "O2" = 0.4377217086785624
    "P2" = 0.4377217086785624
    "Q2" = 5.34136670728
    "R2" = 48.07230036551999
    "S2" = 0.01420895420087014
    "T2" = 0.01420895420087014

    "G3" = 1.328826
    "H3" = 3.986478
    "I3" = 0.03246115949735631
    "J3" = 0.03246115949735631
    "O3" = 0.02575118419467902
    "P3" = 0.02575118419467902
    "Q3" = 0.314232799524
    "R3" = 2.828095195716
    "S3" = 0.0008359132973892763
    "T3" = 0.0008359132973892763

    "G4" = 1.328826
    "H4" = 3.986478
    "I4" = 0.03246115949735631
    "J4" = 0.03246115949735631
    "M4" = 4.926946666666667
    "N4" = 14.78084
    "O4" = 0.5365271071267587
    "P4" = 0.5365271071267587
    "Q4" = 6.547054831280001
    "R4" = 58.92349348152
    "S4" = 0.01741629199909689
    "T4" = 0.01741629199909689

    "I5" = 0.9174921121349238
    "J5" = 0.9174921121349238
    "M5" = 4.019613333333333
    "N5" = 12.05884
    "O5" = 0.4377217086785624
    "P5" = 0.4377217086785624
    "Q5" = 150.9700176405778
    "R5" = 1358.7301587652
    "S5" = 0.401606215022802
    "T5" = 0.401606215022802

    "I6" = 0.9174921121349238
    "J6" = 0.9174921121349238
    "O6" = 0.02575118419467902
    "P6" = 0.02575118419467902
    "S6" = 0.02362650837675152
    "T6" = 0.02362650837675152

    "I7" = 0.9174921121349238
    "J7" = 0.9174921121349238
    "M7" = 4.926946666666667
    "N7" = 14.78084
    "O7" = 0.5365271071267587
    "P7" = 0.5365271071267587
    "Q7" = 185.0479544916889
    "R7" = 1665.4315904252
    "S7" = 0.4922593887353703
    "T7" = 0.4922593887353703

    "G8" = 2.048706666666666
    "H8" = 6.14612
    "I8" = 0.05004672836771996
    "J8" = 0.05004672836771996
    "M8" = 4.019613333333333
    "N8" = 12.05884
    "O8" = 0.4377217086785624
    "P8" = 0.4377217086785624
    "Q8" = 8.23500863342222
    "R8" = 74.11507770079999
    "S8" = 0.02190653945489026
    "T8" = 0.02190653945489026

    "G9" = 2.048706666666666
    "H9" = 6.14612
    "I9" = 0.05004672836771996
    "J9" = 0.05004672836771996
    "O9" = 0.02575118419467902
    "P9" = 0.02575118419467902
    "Q9" = 0.4844658602933333
    "R9" = 4.36019274264
    "S9" = 0.001288762520538224
    "T9" = 0.001288762520538224

    "G10" = 2.048706666666666
    "H10" = 6.14612
    "I10" = 0.05004672836771996
    "J10" = 0.05004672836771996
    "M10" = 4.926946666666667
    "N10" = 14.78084
    "O10" = 0.5365271071267587
    "P10" = 0.5365271071267587
    "Q10" = 10.09386848231111
    "R10" = 90.84481634080001
    "S10" = 0.02685142639229148
    "T10" = 0.02685142639229148
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
